# Apply weekly update: insert a new data row at row 10 (shifting existing
# rows 10-18 down to 11-19), and populate the new row 10 with the latest
# weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 10; this shifts rows 10..18 down to 11..19
# and duplicates the formatting (including the date style) of the row
# above it, consistent with the existing rows.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with this week's data.
$ws.Cells.Item(10, 1).Value = 3
$ws.Cells.Item(10, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(10, 3).Value = "Coquimbo"
$ws.Cells.Item(10, 4).Value = 45001
$ws.Cells.Item(10, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(10, 5).Value = 5
$ws.Cells.Item(10, 6).Value = "Fruta"
$ws.Cells.Item(10, 7).Value = 100101
$ws.Cells.Item(10, 8).Value = "Berries"
$ws.Cells.Item(10, 9).Value = 100101004
$ws.Cells.Item(10, 10).Value = "Frambuesa"
$ws.Cells.Item(10, 11).Value = "Sin especificar"
$ws.Cells.Item(10, 12).Value = "Primera"
$ws.Cells.Item(10, 13).Value = 66
$ws.Cells.Item(10, 14).Value = 7500
$ws.Cells.Item(10, 15).Value = 8000
$ws.Cells.Item(10, 16).Value = 7773
$ws.Cells.Item(10, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(10, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(10, 19).Value = 3886
$ws.Cells.Item(10, 20).Value = 2
